$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column G (old G..O shift right to H..P)
$ws.Columns("G:G").Insert()

# Header for the new "Natural Gas Price (USD/kg)" column
$ws.Range("G1").Value = "Natural Gas Price (USD/kg)"

# New values for rows 2-51 (one per state)
$values = @{
  2=0.19; 3=0.26; 4=0.21; 5=0.34; 6=0.41; 7=0.27; 8=0.34; 9=0.58; 10=0.3;
  11=0.21; 12=1.11; 13=0.2; 14=0.23; 15=0.29; 16=0.24; 17=0.18; 18=0.19; 19=0.14;
  20=0.42; 21=0.48; 22=0.5; 23=0.32; 24=0.22; 25=0.22; 26=0.31; 27=0.32; 28=0.21;
  29=0.32; 30=0.46; 31=0.41; 32=0.17; 33=0.38; 34=0.28; 35=0.11; 36=0.32; 37=0.12;
  38=0.26; 39=0.43; 40=0.49; 41=0.21; 42=0.24; 43=0.22; 44=0.12; 45=0.27; 46=0.19;
  47=0.2; 48=0.42; 49=0.13; 50=0.24; 51=0.22
}

# The data cells use the plain General format (unlike the neighboring columns)
$ws.Range("G2:G51").Style = "Normal"

foreach ($row in $values.Keys) {
  $ws.Cells.Item($row, 7).Value = $values[$row]
}

# Row 53 is a blank spacer row; the inserted column leaves no cell there
$ws.Cells.Item(53, 7).Clear() | Out-Null

# Width of the freshly inserted column
$ws.Columns("G:G").ColumnWidth = 23.3

# Restore the remembered Data > Sort state over the shifted helper columns
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("P3:P26")) | Out-Null
$ws.Sort.SetRange($ws.Range("N3:P26"))
$ws.Sort.Apply()

# Selection left by the author after the edit
$ws.Range("G51").Select() | Out-Null
